$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 33
$ws.Range("F5").Value = 5322
$ws.Range("F6").Value = 5322
$ws.Range("F7").Value = 167
$ws.Range("F11").Value = 1192
$ws.Range("F12").Value = 752
$ws.Range("F13").Value = 5276
$ws.Range("F14").Value = 34
$ws.Range("F15").Value = 77
$ws.Range("F16").Value = 98
$ws.Range("F17").Value = 2588
$ws.Range("F18").Value = 2588
$ws.Range("F19").Value = 261
$ws.Range("F20").Value = 104
$ws.Range("F21").Value = 253
$ws.Range("F22").Value = 3979
$ws.Range("F26").Value = 3900
$ws.Range("F28").Value = 183
$ws.Range("F29").Value = 251
$ws.Range("F30").Value = 214
$ws.Range("F31").Value = 113
$ws.Range("F36").Value = 26
$ws.Range("F37").Value = 6919
$ws.Range("F38").Value = 1131
$ws.Range("F39").Value = 538
$ws.Range("F42").Value = 1407
$ws.Range("F43").Value = 180
$ws.Range("F44").Value = 726
$ws.Range("F46").Value = 2359
$ws.Range("F49").Value = 11
$ws.Range("F50").Value = 790
$ws.Range("F51").Value = 941

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 91
$ws.Range("F16").Value = 10
$ws.Range("F17").Value = 145
$ws.Range("F22").Value = 54
$ws.Range("F25").Value = 826

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 33
$ws.Range("F7").Value = 5322
$ws.Range("F8").Value = 5322
$ws.Range("F9").Value = 167
$ws.Range("F11").Value = 91
$ws.Range("F14").Value = 1192
$ws.Range("F15").Value = 752
$ws.Range("F16").Value = 34
$ws.Range("F17").Value = 77
$ws.Range("F18").Value = 98
$ws.Range("F19").Value = 2588
$ws.Range("F20").Value = 2588
$ws.Range("F21").Value = 261
$ws.Range("F22").Value = 104
$ws.Range("F23").Value = 253
$ws.Range("F24").Value = 3980
$ws.Range("F25").Value = 3900
$ws.Range("F27").Value = 183
$ws.Range("F28").Value = 251
$ws.Range("F29").Value = 214
$ws.Range("F30").Value = 113
$ws.Range("F34").Value = 26
$ws.Range("F35").Value = 145
$ws.Range("F36").Value = 6919
$ws.Range("F37").Value = 1131
$ws.Range("F38").Value = 538
$ws.Range("F42").Value = 1407
$ws.Range("F43").Value = 180
$ws.Range("F44").Value = 726
$ws.Range("F46").Value = 2359
$ws.Range("F49").Value = 790
$ws.Range("F50").Value = 941
